$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new odds value, derived from the authoritative diff.
$updates = @{
    "AC5" = 67
    "AH5" = 41
    "AJ5" = 401
    "J5" = 1.07
    "K5" = 9
    "L5" = 1.4
    "M5" = 2.75
    "N5" = 2.25
    "O5" = 1.62
    "P5" = 1.5
    "Q5" = 2.5
    "R5" = 2
    "S5" = 1.75
    "T5" = 6
    "U5" = 8
    "Y5" = 34
    "Z5" = 8
    "I9" = 4.75
    "N9" = 1.57
    "O9" = 2.35
    "P9" = 1.29
    "Q9" = 3.5
    "R9" = 1.57
    "S9" = 2.25
    "AB10" = 13
    "AD10" = 9.5
    "AE10" = 13
    "AF10" = 10
    "AG10" = 26
    "AH10" = 21
    "G10" = 2.55
    "I10" = 2.55
    "L10" = 1.25
    "M10" = 3.75
    "N10" = 1.85
    "O10" = 1.95
    "R10" = 1.67
    "S10" = 2.1
    "T10" = 9.5
    "U10" = 13
    "V10" = 10
    "W10" = 26
    "X10" = 21
    "Y10" = 29
    "AA11" = 7
    "AE11" = 23
    "G11" = 1.67
    "H11" = 3.7
    "I11" = 5
    "K11" = 9
    "L11" = 1.4
    "M11" = 2.75
    "N11" = 2.2
    "O11" = 1.65
    "P11" = 1.5
    "Q11" = 2.5
    "T11" = 5.5
    "Z11" = 8
    "AB12" = 17
    "AI12" = 41
    "G12" = 1.9
    "H12" = 3.2
    "I12" = 4.33
    "K12" = 7.5
    "N12" = 2.3
    "O12" = 1.6
    "R12" = 2
    "S12" = 1.73
    "G13" = 2.5
    "I13" = 2.8
    "K13" = 9
    "Z13" = 9
    "N14" = 2.15
    "O14" = 1.67
    "AE15" = 26
    "G15" = 1.7
    "H15" = 3.3
    "N15" = 2.6
    "O15" = 1.48
    "X15" = 19
    "J16" = 1.08
    "K16" = 8
    "AD18" = 11
    "AE18" = 17
    "G18" = 2.25
    "I18" = 3.3
    "Y18" = 26
    "AA44" = 7
    "AE44" = 9
    "AG44" = 15
    "G44" = 4
    "H44" = 3.6
    "I44" = 1.9
    "L44" = 1.25
    "M44" = 3.75
    "N44" = 1.9
    "O44" = 1.95
    "X44" = 34
    "AC45" = 32
    "AE45" = 16
    "AF45" = 9.5
    "AG45" = 29
    "AH45" = 17.5
    "AI45" = 19.5
    "G45" = 2.55
    "I45" = 2.45
    "U45" = 17.5
    "V45" = 9.75
    "W45" = 32
    "X45" = 18
    "Y45" = 19.5
    "AG46" = 41
    "I46" = 3.25
    "J46" = 1.08
    "K46" = 8
    "W46" = 19
    "X46" = 19
    "AE47" = 13
    "AH47" = 29
    "G47" = 2.4
    "I47" = 2.88
    "K47" = 6.2
    "Q47" = 2.5
    "R47" = 1.93
    "S47" = 1.78
    "T47" = 6.5
    "V47" = 10
    "W47" = 23
    "Z47" = 7
    "N48" = 1.65
    "O48" = 2.2
    "P48" = 1.33
    "P49" = 1.4
    "AA70" = 5.8
    "AB70" = 15
    "AC70" = 80
    "AD70" = 8.75
    "AE70" = 19
    "AF70" = 13
    "AG70" = 60
    "AH70" = 40
    "AJ70" = 800
    "G70" = 2.07
    "H70" = 2.95
    "I70" = 3.7
    "L70" = 1.39
    "M70" = 2.57
    "N70" = 2.12
    "O70" = 1.57
    "R70" = 1.83
    "S70" = 1.78
    "T70" = 6.5
    "U70" = 9.5
    "V70" = 8.5
    "W70" = 20
    "X70" = 17.5
    "Y70" = 30
    "AA71" = 7.1
    "AB71" = 21
    "AD71" = 12
    "G71" = 1.6
    "H71" = 3.55
    "I71" = 5.5
    "L71" = 1.37
    "M71" = 2.62
    "N71" = 2.07
    "O71" = 1.6
    "P71" = 1.42
    "Q71" = 2.45
    "R71" = 2.07
    "S71" = 1.6
    "T71" = 5.4
    "U71" = 6.4
    "V71" = 8.5
    "W71" = 11.25
    "Y71" = 37
    "Z71" = 8
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
